$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for this product; insert it at the top
# of the data block (row 113), pushing every existing record down by one
# row (old row 113 -> 114, ..., old row 136 -> 137).
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new record's data.
$ws.Cells.Item(113, 1).Value = 1
$ws.Cells.Item(113, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(113, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(113, 4).Value = 45124
$ws.Cells.Item(113, 5).Value = 15
$ws.Cells.Item(113, 6).Value = 100112038
$ws.Cells.Item(113, 7).Value = "Cebollín baby"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 250
$ws.Cells.Item(113, 11).Value = 1400
$ws.Cells.Item(113, 12).Value = 1500
$ws.Cells.Item(113, 13).Value = 1440
$ws.Cells.Item(113, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(113, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(113, 16).Value = 720
$ws.Cells.Item(113, 17).Value = 2
$ws.Cells.Item(113, 18).Value = "Hortaliza"
